# Update titles and descriptions on the "Common" sheet.
#
# Splits the old "KVM RAM" section (rows 85-92) into two sections:
#   "KVM and VCenter RAM" (VSD/VSC/VSTAT RAM - valid for KVM+VCenter)
#   "KVM RAM" (VCIN/NUH/Webfilter/Portal RAM - valid for KVM only)
# and the old "CPU" section (rows 93-101) into two sections:
#   "KVM and VCenter CPU" (VSD/VSC/VSTAT/VNSUTIL CPU cores - valid for KVM+VCenter)
#   "KVM CPU" (NUH/VCIN/Portal/Webfilter CPU cores - valid for KVM only)
# This requires inserting two new header rows, and renaming several row
# labels and comment bodies.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# --- Insert the two new section-header rows -------------------------------
# Insert in ascending order: the first insert (at row 89) shifts the target
# location of the second header (originally row 98/99) down to row 99, so it
# must be applied before the second insert.

# New row 89: "KVM RAM" header (splits the old RAM block)
$ws.Rows.Item(89).Insert()
$ws.Range("A85").Copy()
$ws.Range("A89").PasteSpecial(-4122)
$ws.Range("A89:B89").Merge()
$ws.Range("B89").ClearContents()
$ws.Range("A89").Value = "KVM RAM"

# New row 99: "KVM CPU" header (splits the old CPU block)
$ws.Rows.Item(99).Insert()
$ws.Range("A94").Copy()
$ws.Range("A99").PasteSpecial(-4122)
$ws.Range("A99:B99").Merge()
$ws.Range("B99").ClearContents()
$ws.Range("A99").Value = "KVM CPU"

# --- Rename section headers and row labels ---------------------------------
$ws.Range("A85").Value = "KVM and VCenter RAM"
$ws.Range("A86").Value = "KVM VSD RAM"
$ws.Range("A87").Value = "KVM VSC RAM"
$ws.Range("A88").Value = "KVM VSTAT RAM"

$ws.Range("A94").Value = "KVM and VCenter CPU"
$ws.Range("A95").Value = "KVM VSD CPU cores"
$ws.Range("A96").Value = "KVM VSC CPU cores"
$ws.Range("A97").Value = "KVM VSTAT CPU cores"
$ws.Range("A98").Value = "KVM VNSUTIL CPU cores"

# --- Update comment text on the renamed RAM rows ----------------------------
$ws.Range("A86").Comment.Text("For KVM and VCenter deployments: amount of VSD RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]")
$ws.Range("A87").Comment.Text("For KVM and VCenter deployments: amount of VSC RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 4]")
$ws.Range("A88").Comment.Text("For KVM and VCenter deployments: amount of VSTAT RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 16]")

# --- Update comment text on the renamed CPU rows ----------------------------
$ws.Range("A95").Comment.Text("For KVM and VCenter deployments: number of CPU's for VSD. [default: 6]")
$ws.Range("A96").Comment.Text("For KVM and VCenter deployments: number of CPU's for VSC. [default: 6]")
$ws.Range("A97").Comment.Text("For KVM and VCenter deployments: number of CPU's for VSTAT. [default: 6]")
$ws.Range("A98").Comment.Text("For KVM and VCenter deployments: number of CPU's for VNSUTIL. [default: 2]")

Write-Output "done"
